# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G) is regenerated from the updated source data.
# Only the G2:G62 values change; everything else on the sheet is left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @(1,1,2,2,0,1,0,0,3,1,2,0,2,1,0,1,1,3,3,1,4,2,0,1,0,0,1,0,1,1,0,0,0,2,0,1,2,1,1,2,0,1,0,1,0,0,2,1,0,0,1,1,1,0,0,0,2,2,1,0,2)

for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
